$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24050
$ws.Range("D2").Value = -0.0494
$ws.Range("I2").Value = 4.16
$ws.Range("J2").Value = 84
$ws.Range("K2").Value = 84

$ws.Range("C3").Value = 104100
$ws.Range("D3").Value = -0.0086
$ws.Range("I3").Value = 6.24
$ws.Range("J3").Value = 68
$ws.Range("K3").Value = 68

$ws.Range("C4").Value = 456000
$ws.Range("I4").Value = 4.17
$ws.Range("J4").Value = 81
$ws.Range("K4").Value = 81

$ws.Range("C5").Value = 33550
$ws.Range("D5").Value = -0.0015
$ws.Range("I5").Value = 5.96
$ws.Range("J5").Value = 54
$ws.Range("K5").Value = 54

$ws.Range("C6").Value = 33400
$ws.Range("D6").Value = -0.0162
$ws.Range("I6").Value = 3.59
$ws.Range("J6").Value = 89
$ws.Range("K6").Value = 89

$ws.Range("C7").Value = 27150
$ws.Range("D7").Value = -0.0037
$ws.Range("I7").Value = 4.42
$ws.Range("J7").Value = 78
$ws.Range("K7").Value = 78

$ws.Range("C8").Value = 11030
$ws.Range("D8").Value = 0.0046
$ws.Range("I8").Value = 4.67

$ws.Range("C9").Value = 76800
$ws.Range("D9").Value = -0.0192
$ws.Range("I9").Value = 3.91
$ws.Range("J9").Value = 63
$ws.Range("K9").Value = 63

$ws.Range("C10").Value = 216500
$ws.Range("D10").Value = -0.0046
$ws.Range("I10").Value = 5.54
$ws.Range("J10").Value = 53
$ws.Range("K10").Value = 53

$ws.Range("C11").Value = 134200
$ws.Range("D11").Value = 0.0167
$ws.Range("I11").Value = 5.07

$ws.Range("C12").Value = 20700
$ws.Range("D12").Value = 0.0049
$ws.Range("I12").Value = 4.59
$ws.Range("J12").Value = 82
$ws.Range("K12").Value = 82

$ws.Range("C13").Value = 72400
$ws.Range("D13").Value = 0.0112

$ws.Range("C14").Value = 55500
$ws.Range("D14").Value = 0
$ws.Range("D14").NumberFormat = "0%"
$ws.Range("I14").Value = 6.38
$ws.Range("J14").Value = 71
$ws.Range("K14").Value = 71

$ws.Range("C15").Value = 81900
$ws.Range("D15").Value = -0.0012
$ws.Range("I15").Value = 6.72
$ws.Range("J15").Value = 86
$ws.Range("K15").Value = 86

$ws.Range("C16").Value = 19850
$ws.Range("D16").Value = -0.0198
$ws.Range("I16").Value = 5.37
$ws.Range("J16").Value = 83
$ws.Range("K16").Value = 83

$ws.Range("C17").Value = 52600
$ws.Range("D17").Value = -0.0131
$ws.Range("I17").Value = 5.32
$ws.Range("J17").Value = 78
$ws.Range("K17").Value = 78

$ws.Range("C18").Value = 21150
$ws.Range("D18").Value = -0.0047
$ws.Range("D18").NumberFormat = "0.00%"
$ws.Range("I18").Value = 5.82
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = 46

$ws.Range("C19").Value = 55900
$ws.Range("D19").Value = -0.0089
$ws.Range("I19").Value = 3.58
$ws.Range("J19").Value = 91
$ws.Range("K19").Value = 91

$ws.Range("C20").Value = 14670
$ws.Range("D20").Value = 0.0055
$ws.Range("I20").Value = 4.43
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 77

$ws.Range("C21").Value = 132500
$ws.Range("D21").Value = 0.0076
$ws.Range("I21").Value = 4.08
$ws.Range("J21").Value = 83
$ws.Range("K21").Value = 83

$ws.Range("C22").Value = 45250
$ws.Range("D22").Value = 0.0134
$ws.Range("I22").Value = 3.22
$ws.Range("J22").Value = 54
$ws.Range("K22").Value = 54

$ws.Range("C23").Value = 71200
$ws.Range("D23").Value = 0.0274
$ws.Range("I23").Value = 3.03
$ws.Range("J23").Value = 95
$ws.Range("K23").Value = 95

$ws.Range("C24").Value = 51300
$ws.Range("D24").Value = 0.002
$ws.Range("I24").Value = 5.26
$ws.Range("J24").Value = 77
$ws.Range("K24").Value = 77

$ws.Range("C25").Value = 92500
$ws.Range("D25").Value = 0.0165
$ws.Range("I25").Value = 3.89
$ws.Range("J25").Value = 93
$ws.Range("K25").Value = 93

$ws.Range("C26").Value = 118800
$ws.Range("D26").Value = 0.0137
$ws.Range("I26").Value = 2.67
$ws.Range("J26").Value = 91
$ws.Range("K26").Value = 91

$ws.Range("C27").Value = 15180
$ws.Range("D27").Value = -0.0181
$ws.Range("I27").Value = 4.28
$ws.Range("J27").Value = 92
$ws.Range("K27").Value = 92

$ws.Range("C28").Value = 14540
$ws.Range("D28").Value = 0.0048
$ws.Range("I28").Value = 3.44
$ws.Range("J28").Value = 90
$ws.Range("K28").Value = 90

$ws.Range("C29").Value = 23750
$ws.Range("D29").Value = -0.0246
$ws.Range("I29").Value = 4.19
$ws.Range("J29").Value = 89
$ws.Range("K29").Value = 89

$ws.Range("C30").Value = 25550
$ws.Range("D30").Value = 0.0059
$ws.Range("I30").Value = 4.7
$ws.Range("J30").Value = 92
$ws.Range("K30").Value = 92

$ws.Range("A1").Select()
